# Update crypto price/volume data per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.659.89'
$ws.Range("E2").Value = '  -2.12%  '
$ws.Range("D3").Value = '1.795.79'
$ws.Range("E3").Value = '  -1.91%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''231.58'
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("D6").Value = '''0.5874'
$ws.Range("D7").Value = '''1.005'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.2765'
$ws.Range("E8").Value = '  -0.97%  '
$ws.Range("D9").Value = '''0.06749'
$ws.Range("E9").Value = '  -4.32%  '
$ws.Range("D10").Value = '''23.17'
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("D11").Value = '''0.07523'
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").Value = '1.796.49'
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("D13").Value = '''4.783'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '''0.6151'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").Value = '2.040.32'
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '''75.24'
$ws.Range("E16").Value = '  -4.93%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.000009041'
$ws.Range("E17").Value = '  -8.59%  '
$ws.Range("D18").Value = '28.647.12'
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("D19").Value = '''5.463'
$ws.Range("E19").Value = '  -6.61%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '''209.87'
$ws.Range("E21").Value = '  -6.58%  '
$ws.Range("D22").Value = '''11.48'
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("D23").Value = '''6.815'
$ws.Range("E23").Value = '  -2.93%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '''153.51'
$ws.Range("E25").Value = '  -1.56%  '
$ws.Range("D26").Value = '''8.083'
$ws.Range("E26").Value = '  +1.21%  '
$ws.Range("D27").Value = '''0.1258'
$ws.Range("D28").Value = '''16.41'
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("E29").Value = '  -3.69%  '
$ws.Range("D30").Value = '''0.06106'
$ws.Range("E30").Value = '  -4.31%  '
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("D32").Value = '''3.806'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").Value = '''3.784'
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").Value = '''1.731'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").Value = '''1.045'
$ws.Range("E35").Value = '  -5.85%  '
$ws.Range("D36").Value = '''0.6402'
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("E37").Value = '  -1.85%  '
$ws.Range("D38").Value = '''2.715'
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").Value = '''6.412'
$ws.Range("E39").Value = '  -2.01%  '
$ws.Range("D40").Value = '''0.01697'
$ws.Range("E40").Value = '  -3.23%  '
$ws.Range("D41").Value = '1.142.40'
$ws.Range("E41").Value = '  -6.29%  '
$ws.Range("D42").Value = '''0.8794'
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("D43").Value = '''1.007'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '''100.05'
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("D45").Value = '1.947.94'
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("D46").Value = '''60.00'
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("D47").Value = '''0.00000000112'
$ws.Range("E47").Value = '  -3.64%  '
$ws.Range("D48").Value = '''1.583'
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.05490'
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''8.336'
$ws.Range("E50").Value = '  -2.61%  '
$ws.Range("D51").Value = '''0.4484'
$ws.Range("E51").Value = '  -1.58%  '
